# Weekly driver report update for 2025-04-19
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: driver version bump + updated counts/percentage
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.20.1.1"
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 98.90000000000001

# Row 4 totals: Critical Minutes total follows row 3's drop
$ws.Range("C4").Value = 3

# Row 12: clear the stale Driver Vintage date
$ws.Range("E12").ClearContents()

# Row 14: updated Total Samples count
$ws.Range("B14").Value = 265400
